$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1 header text is unchanged ("running_distance") - only its underlying shared
# string index shifts because unused English enum strings get removed below.

# Row 2: English option codes replaced by their Hebrew display values
$ws.Range("E2").Value = "ארוך"
$ws.Range("F2").Value = "ריצה"
$ws.Range("H2").Value = "כביש"

# Row 3: English option codes replaced by their Hebrew display values
$ws.Range("E3").Value = "בינוני"
$ws.Range("F3").Value = "לכל היום"
$ws.Range("H3").Value = "כביש"

# New RFID-related columns AX (comp) / AY (damping)
$ws.Range("AX1").Value = "comp"
$ws.Range("AY1").Value = "damping"

$ws.Range("AX2").Value = "לא"
$ws.Range("AY2").Value = "קל"

$ws.Range("AX3").Value = "כן"
$ws.Range("AY3").Value = "מקסימלי"

# Match the author's final selection in the saved view
$ws.Range("AX3").Select() | Out-Null
